# Update the Efna5/Epha7 LR-pair sheet with refreshed TPM-derived values.
# The sending/target cluster combinations now form a full 3x3 grid
# (ECs, FAPs, MuSCs) instead of the previous partial 2x3 layout, so the
# sheet grows from 7 to 10 used rows (A1:T10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = 'ECs'  # A2
$ws.Cells.Item(2, 2).Value2 = 'Efna5'  # B2
$ws.Cells.Item(2, 3).Value2 = 'Epha7'  # C2
$ws.Cells.Item(2, 4).Value2 = 'ECs'  # D2
$ws.Cells.Item(2, 5).Value2 = 1  # E2
$ws.Cells.Item(2, 6).Value2 = 0.3333333333333333  # F2
$ws.Cells.Item(2, 7).Value2 = 0.05800433333333333  # G2
$ws.Cells.Item(2, 8).Value2 = 0.174013  # H2
$ws.Cells.Item(2, 9).Value2 = 0.02087975181349295  # I2
$ws.Cells.Item(2, 10).Value2 = 0.02087975181349295  # J2
$ws.Cells.Item(2, 11).Value2 = 3  # K2
$ws.Cells.Item(2, 12).Value2 = 1  # L2
$ws.Cells.Item(2, 13).Value2 = 0.1944653333333334  # M2
$ws.Cells.Item(2, 14).Value2 = 0.583396  # N2
$ws.Cells.Item(2, 15).Value2 = 0.04942840076761122  # O2
$ws.Cells.Item(2, 16).Value2 = 0.04942840076761121  # P2
$ws.Cells.Item(2, 17).Value2 = 0.01127983201644444  # Q2
$ws.Cells.Item(2, 18).Value2 = 0.101518488148  # R2
$ws.Cells.Item(2, 19).Value2 = 0.001032052740565587  # S2
$ws.Cells.Item(2, 20).Value2 = 0.001032052740565587  # T2

# Row 3
$ws.Cells.Item(3, 1).Value2 = 'ECs'  # A3
$ws.Cells.Item(3, 2).Value2 = 'Efna5'  # B3
$ws.Cells.Item(3, 3).Value2 = 'Epha7'  # C3
$ws.Cells.Item(3, 4).Value2 = 'FAPs'  # D3
$ws.Cells.Item(3, 5).Value2 = 1  # E3
$ws.Cells.Item(3, 6).Value2 = 0.3333333333333333  # F3
$ws.Cells.Item(3, 7).Value2 = 0.05800433333333333  # G3
$ws.Cells.Item(3, 8).Value2 = 0.174013  # H3
$ws.Cells.Item(3, 9).Value2 = 0.02087975181349295  # I3
$ws.Cells.Item(3, 10).Value2 = 0.02087975181349295  # J3
$ws.Cells.Item(3, 11).Value2 = 3  # K3
$ws.Cells.Item(3, 12).Value2 = 1  # L3
$ws.Cells.Item(3, 13).Value2 = 0.1693566666666667  # M3
$ws.Cells.Item(3, 14).Value2 = 0.50807  # N3
$ws.Cells.Item(3, 15).Value2 = 0.04304638286515546  # O3
$ws.Cells.Item(3, 16).Value2 = 0.04304638286515546  # P3
$ws.Cells.Item(3, 17).Value2 = 0.009823420545555557  # Q3
$ws.Cells.Item(3, 18).Value2 = 0.08841078491  # R3
$ws.Cells.Item(3, 19).Value2 = 0.0008987977906930416  # S3
$ws.Cells.Item(3, 20).Value2 = 0.0008987977906930415  # T3

# Row 4
$ws.Cells.Item(4, 1).Value2 = 'ECs'  # A4
$ws.Cells.Item(4, 2).Value2 = 'Efna5'  # B4
$ws.Cells.Item(4, 3).Value2 = 'Epha7'  # C4
$ws.Cells.Item(4, 4).Value2 = 'MuSCs'  # D4
$ws.Cells.Item(4, 5).Value2 = 1  # E4
$ws.Cells.Item(4, 6).Value2 = 0.3333333333333333  # F4
$ws.Cells.Item(4, 7).Value2 = 0.05800433333333333  # G4
$ws.Cells.Item(4, 8).Value2 = 0.174013  # H4
$ws.Cells.Item(4, 9).Value2 = 0.02087975181349295  # I4
$ws.Cells.Item(4, 10).Value2 = 0.02087975181349295  # J4
$ws.Cells.Item(4, 11).Value2 = 3  # K4
$ws.Cells.Item(4, 12).Value2 = 1  # L4
$ws.Cells.Item(4, 13).Value2 = 3.570461333333333  # M4
$ws.Cells.Item(4, 14).Value2 = 10.711384  # N4
$ws.Cells.Item(4, 15).Value2 = 0.9075252163672334  # O4
$ws.Cells.Item(4, 16).Value2 = 0.9075252163672333  # P4
$ws.Cells.Item(4, 17).Value2 = 0.2071022293324444  # Q4
$ws.Cells.Item(4, 18).Value2 = 1.863920063992  # R4
$ws.Cells.Item(4, 19).Value2 = 0.01894890128223432  # S4
$ws.Cells.Item(4, 20).Value2 = 0.01894890128223432  # T4

# Row 5
$ws.Cells.Item(5, 1).Value2 = 'FAPs'  # A5
$ws.Cells.Item(5, 2).Value2 = 'Efna5'  # B5
$ws.Cells.Item(5, 3).Value2 = 'Epha7'  # C5
$ws.Cells.Item(5, 4).Value2 = 'ECs'  # D5
$ws.Cells.Item(5, 5).Value2 = 3  # E5
$ws.Cells.Item(5, 6).Value2 = 1  # F5
$ws.Cells.Item(5, 7).Value2 = 1.666083666666667  # G5
$ws.Cells.Item(5, 8).Value2 = 4.998251  # H5
$ws.Cells.Item(5, 9).Value2 = 0.5997381826733804  # I5
$ws.Cells.Item(5, 10).Value2 = 0.5997381826733805  # J5
$ws.Cells.Item(5, 11).Value2 = 3  # K5
$ws.Cells.Item(5, 12).Value2 = 1  # L5
$ws.Cells.Item(5, 13).Value2 = 0.1944653333333334  # M5
$ws.Cells.Item(5, 14).Value2 = 0.583396  # N5
$ws.Cells.Item(5, 15).Value2 = 0.04942840076761122  # O5
$ws.Cells.Item(5, 16).Value2 = 0.04942840076761121  # P5
$ws.Cells.Item(5, 17).Value2 = 0.3239955155995556  # Q5
$ws.Cells.Item(5, 18).Value2 = 2.915959640396  # R5
$ws.Cells.Item(5, 19).Value2 = 0.02964409924881867  # S5
$ws.Cells.Item(5, 20).Value2 = 0.02964409924881868  # T5

# Row 6
$ws.Cells.Item(6, 1).Value2 = 'FAPs'  # A6
$ws.Cells.Item(6, 2).Value2 = 'Efna5'  # B6
$ws.Cells.Item(6, 3).Value2 = 'Epha7'  # C6
$ws.Cells.Item(6, 4).Value2 = 'FAPs'  # D6
$ws.Cells.Item(6, 5).Value2 = 3  # E6
$ws.Cells.Item(6, 6).Value2 = 1  # F6
$ws.Cells.Item(6, 7).Value2 = 1.666083666666667  # G6
$ws.Cells.Item(6, 8).Value2 = 4.998251  # H6
$ws.Cells.Item(6, 9).Value2 = 0.5997381826733804  # I6
$ws.Cells.Item(6, 10).Value2 = 0.5997381826733805  # J6
$ws.Cells.Item(6, 11).Value2 = 3  # K6
$ws.Cells.Item(6, 12).Value2 = 1  # L6
$ws.Cells.Item(6, 13).Value2 = 0.1693566666666667  # M6
$ws.Cells.Item(6, 14).Value2 = 0.50807  # N6
$ws.Cells.Item(6, 15).Value2 = 0.04304638286515546  # O6
$ws.Cells.Item(6, 16).Value2 = 0.04304638286515546  # P6
$ws.Cells.Item(6, 17).Value2 = 0.2821623761744445  # Q6
$ws.Cells.Item(6, 18).Value2 = 2.53946138557  # R6
$ws.Cells.Item(6, 19).Value2 = 0.02581655943021088  # S6
$ws.Cells.Item(6, 20).Value2 = 0.02581655943021088  # T6

# Row 7
$ws.Cells.Item(7, 1).Value2 = 'FAPs'  # A7
$ws.Cells.Item(7, 2).Value2 = 'Efna5'  # B7
$ws.Cells.Item(7, 3).Value2 = 'Epha7'  # C7
$ws.Cells.Item(7, 4).Value2 = 'MuSCs'  # D7
$ws.Cells.Item(7, 5).Value2 = 3  # E7
$ws.Cells.Item(7, 6).Value2 = 1  # F7
$ws.Cells.Item(7, 7).Value2 = 1.666083666666667  # G7
$ws.Cells.Item(7, 8).Value2 = 4.998251  # H7
$ws.Cells.Item(7, 9).Value2 = 0.5997381826733804  # I7
$ws.Cells.Item(7, 10).Value2 = 0.5997381826733805  # J7
$ws.Cells.Item(7, 11).Value2 = 3  # K7
$ws.Cells.Item(7, 12).Value2 = 1  # L7
$ws.Cells.Item(7, 13).Value2 = 3.570461333333333  # M7
$ws.Cells.Item(7, 14).Value2 = 10.711384  # N7
$ws.Cells.Item(7, 15).Value2 = 0.9075252163672334  # O7
$ws.Cells.Item(7, 16).Value2 = 0.9075252163672333  # P7
$ws.Cells.Item(7, 17).Value2 = 5.948687309931554  # Q7
$ws.Cells.Item(7, 18).Value2 = 53.53818578938399  # R7
$ws.Cells.Item(7, 19).Value2 = 0.5442775239943509  # S7
$ws.Cells.Item(7, 20).Value2 = 0.5442775239943509  # T7

# Row 8
$ws.Cells.Item(8, 1).Value2 = 'MuSCs'  # A8
$ws.Cells.Item(8, 2).Value2 = 'Efna5'  # B8
$ws.Cells.Item(8, 3).Value2 = 'Epha7'  # C8
$ws.Cells.Item(8, 4).Value2 = 'ECs'  # D8
$ws.Cells.Item(8, 5).Value2 = 3  # E8
$ws.Cells.Item(8, 6).Value2 = 1  # F8
$ws.Cells.Item(8, 7).Value2 = 1.053930333333333  # G8
$ws.Cells.Item(8, 8).Value2 = 3.161791  # H8
$ws.Cells.Item(8, 9).Value2 = 0.3793820655131266  # I8
$ws.Cells.Item(8, 10).Value2 = 0.3793820655131266  # J8
$ws.Cells.Item(8, 11).Value2 = 3  # K8
$ws.Cells.Item(8, 12).Value2 = 1  # L8
$ws.Cells.Item(8, 13).Value2 = 0.1944653333333334  # M8
$ws.Cells.Item(8, 14).Value2 = 0.583396  # N8
$ws.Cells.Item(8, 15).Value2 = 0.04942840076761122  # O8
$ws.Cells.Item(8, 16).Value2 = 0.04942840076761121  # P8
$ws.Cells.Item(8, 17).Value2 = 0.2049529135817778  # Q8
$ws.Cells.Item(8, 18).Value2 = 1.844576222236  # R8
$ws.Cells.Item(8, 19).Value2 = 0.01875224877822696  # S8
$ws.Cells.Item(8, 20).Value2 = 0.01875224877822695  # T8

# Row 9
$ws.Cells.Item(9, 1).Value2 = 'MuSCs'  # A9
$ws.Cells.Item(9, 2).Value2 = 'Efna5'  # B9
$ws.Cells.Item(9, 3).Value2 = 'Epha7'  # C9
$ws.Cells.Item(9, 4).Value2 = 'FAPs'  # D9
$ws.Cells.Item(9, 5).Value2 = 3  # E9
$ws.Cells.Item(9, 6).Value2 = 1  # F9
$ws.Cells.Item(9, 7).Value2 = 1.053930333333333  # G9
$ws.Cells.Item(9, 8).Value2 = 3.161791  # H9
$ws.Cells.Item(9, 9).Value2 = 0.3793820655131266  # I9
$ws.Cells.Item(9, 10).Value2 = 0.3793820655131266  # J9
$ws.Cells.Item(9, 11).Value2 = 3  # K9
$ws.Cells.Item(9, 12).Value2 = 1  # L9
$ws.Cells.Item(9, 13).Value2 = 0.1693566666666667  # M9
$ws.Cells.Item(9, 14).Value2 = 0.50807  # N9
$ws.Cells.Item(9, 15).Value2 = 0.04304638286515546  # O9
$ws.Cells.Item(9, 16).Value2 = 0.04304638286515546  # P9
$ws.Cells.Item(9, 17).Value2 = 0.1784901281522223  # Q9
$ws.Cells.Item(9, 18).Value2 = 1.60641115337  # R9
$ws.Cells.Item(9, 19).Value2 = 0.01633102564425154  # S9
$ws.Cells.Item(9, 20).Value2 = 0.01633102564425154  # T9

# Row 10
$ws.Cells.Item(10, 1).Value2 = 'MuSCs'  # A10
$ws.Cells.Item(10, 2).Value2 = 'Efna5'  # B10
$ws.Cells.Item(10, 3).Value2 = 'Epha7'  # C10
$ws.Cells.Item(10, 4).Value2 = 'MuSCs'  # D10
$ws.Cells.Item(10, 5).Value2 = 3  # E10
$ws.Cells.Item(10, 6).Value2 = 1  # F10
$ws.Cells.Item(10, 7).Value2 = 1.053930333333333  # G10
$ws.Cells.Item(10, 8).Value2 = 3.161791  # H10
$ws.Cells.Item(10, 9).Value2 = 0.3793820655131266  # I10
$ws.Cells.Item(10, 10).Value2 = 0.3793820655131266  # J10
$ws.Cells.Item(10, 11).Value2 = 3  # K10
$ws.Cells.Item(10, 12).Value2 = 1  # L10
$ws.Cells.Item(10, 13).Value2 = 3.570461333333333  # M10
$ws.Cells.Item(10, 14).Value2 = 10.711384  # N10
$ws.Cells.Item(10, 15).Value2 = 0.9075252163672334  # O10
$ws.Cells.Item(10, 16).Value2 = 0.9075252163672333  # P10
$ws.Cells.Item(10, 17).Value2 = 3.763017503193777  # Q10
$ws.Cells.Item(10, 18).Value2 = 33.86715752874399  # R10
$ws.Cells.Item(10, 19).Value2 = 0.3442987910906481  # S10
$ws.Cells.Item(10, 20).Value2 = 0.3442987910906481  # T10
